$wb = $excel.ActiveWorkbook

# The same edit is applied to both the "展览" sheet and the "全部类型" sheet,
# which contain identical data. Only the F24 target value differs between
# the two sheets (95 vs 96) per the source diff.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- 1. Insert a new record row at row 28 -----------------------------
    # Shift existing rows 28..36 down to 29..37 by copying whole row ranges
    # from the bottom up (so we never overwrite data before it is copied).
    for ($r = 36; $r -ge 28; $r--) {
        $src = $ws.Range("A" + $r + ":I" + $r)
        $dst = $ws.Range("A" + ($r + 1) + ":I" + ($r + 1))
        $src.Copy($dst)
    }

    # --- 2. Populate the newly freed row 28 with the new record -----------
    # Column A (running index) keeps the value that was already there
    # (27), matching the rest of the table's "row - 1" numbering.
    $ws.Range("B28").Value = "2024-05-18"
    $ws.Range("C28").Value = "赣州·原铁崩only"
    $ws.Range("D28").Value = "金岭东大道18号 万达广场(赣州经开店)"
    $ws.Range("E28").Value = "2024.05.18 10:00-05.19 17:00"
    $ws.Range("F28").Value = 0
    $ws.Range("G28").Value = 60
    $ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=84721"
    $ws.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202404/0n0MQiZh1713505673648.jpeg"

    # --- 3. Refresh the "want to go" counters for several existing rows ---
    $ws.Range("F3").Value = 733
    $ws.Range("F8").Value = 1663
    $ws.Range("F9").Value = 6096
    $ws.Range("F11").Value = 350
    $ws.Range("F13").Value = 86
    $ws.Range("F16").Value = 5432
    $ws.Range("F17").Value = 260
    $ws.Range("F22").Value = 97
    $ws.Range("F23").Value = 258

    if ($name -eq "展览") {
        $ws.Range("F24").Value = 95
    } else {
        $ws.Range("F24").Value = 96
    }
}
